# Update "想去人数" (want-to-go count) values in column F across sheets.
# Sheet 1 = 展览 (Exhibition)
# Sheet 2 = 演出 (Performance)
# Sheet 3 = 本地生活 (Local life) - unchanged
# Sheet 4 = 全部类型 (All types)

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 1204
$ws1.Range("F3").Value = 647
$ws1.Range("F5").Value = 0
$ws1.Range("F7").Value = 0
$ws1.Range("F10").Value = 87
$ws1.Range("F11").Value = 0
$ws1.Range("F12").Value = 74

$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 0
$ws2.Range("F3").Value = 0
$ws2.Range("F4").Value = 0
$ws2.Range("F5").Value = 7
$ws2.Range("F6").Value = 0

$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F4").Value = 350
$ws4.Range("F5").Value = 14
$ws4.Range("F6").Value = 20
$ws4.Range("F8").Value = 0
$ws4.Range("F9").Value = 8
$ws4.Range("F10").Value = 0
$ws4.Range("F11").Value = 0
$ws4.Range("F12").Value = 529
$ws4.Range("F13").Value = 87
$ws4.Range("F15").Value = 0
$ws4.Range("F16").Value = 670
